$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full target contents for the module listing table (A1:E13).
$data = @(
    @("Code", "Name", "ChefModule", "ElementName1", "ElementName2"),
    @("GIL11", "pede. Suspendisse dui.", "EL Haddad", "Nullam feugiat placerat", "varius et, euismod"),
    @("GIL12", "a nunc. In", "Badir", "sodales nisi magna", "elementum sem, vitae"),
    @("GIL13", "amet metus. Aliquam", "Ezzine", "Cras vulputate velit", "scelerisque neque sed"),
    @("GIL14", "quam vel sapien", "El Alami Hassoun", "Nunc mauris elit,", "libero et tristique"),
    @("GIL15", "feugiat nec, diam.", "Lazaar", "pellentesque. Sed dictum.", "ridiculus mus. Proin"),
    @("GIL16", "nonummy. Fusce fermentum", "El Haddad", "neque pellentesque massa", "Mauris eu turpis."),
    @("GIL21", "a, arcu. Sed", "EL Haddad", "sit amet risus.", "Nulla facilisi. Sed"),
    @("GIL22", "Suspendisse eleifend. Cras", "El Alami Hassoun", "velit dui, semper", "ligula elit, pretium"),
    @("GIL23", "ante. Nunc mauris", "Badir", "tortor at risus.", "felis. Donec tempor,"),
    @("GIL24", "lobortis quam a", "Ezzine", "euismod est arcu", "ligula eu enim."),
    @("GIL25", "rhoncus. Nullam velit", "Ben Achrab", "ut dolor dapibus", "commodo tincidunt nibh."),
    @("GIL26", "Donec tincidunt. Donec", "EL Haddad", "ornare tortor at", "ac, feugiat non,")
)

# Clear the previously used range (old table spanned A1:F13) before writing
# the new, wider table so the stale column F contents are dropped.
$ws.Range("A1:F13").ClearContents()

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowVals[$c]
    }
}

# Column widths to match the new content (bestFit-style widths from the
# saved workbook). The host's ColumnWidth setter quantises the stored
# OOXML width to the nearest 1/6 character, so the inputs below are the
# pre-images that land closest to the recorded widths of 5.7109375,
# 27, 16.140625, 25.140625, 23.5703125 and 14.7109375 respectively.
$ws.Columns.Item(1).ColumnWidth = 4.83333333333333
$ws.Columns.Item(2).ColumnWidth = 26.1666666666667
$ws.Columns.Item(3).ColumnWidth = 15.3333333333333
$ws.Columns.Item(4).ColumnWidth = 24.3333333333333
$ws.Columns.Item(5).ColumnWidth = 22.6666666666667
$ws.Columns.Item(6).ColumnWidth = 13.8333333333333

# Update the active selection as recorded in the saved workbook.
$ws.Range("G9").Select()
